$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised daily stats for existing rows (F = AgTests, G = AgPosit) ---
$ws.Range("F393").Value = 308410
$ws.Range("F394").Value = 166094
$ws.Range("F395").Value = 752767
$ws.Range("F396").Value = 164972
$ws.Range("F397").Value = 107886
$ws.Range("F400").Value = 150496
$ws.Range("F401").Value = 272215
$ws.Range("G401").Value = 933
$ws.Range("F402").Value = 719662
$ws.Range("G402").Value = 1392
$ws.Range("F403").Value = 351736
$ws.Range("F404").Value = 224673
$ws.Range("F405").Value = 173337
$ws.Range("G405").Value = 688
$ws.Range("F406").Value = 170400
$ws.Range("F407").Value = 157402
$ws.Range("G407").Value = 669
$ws.Range("F408").Value = 303923
$ws.Range("F409").Value = 705696
$ws.Range("G409").Value = 1003
$ws.Range("F410").Value = 363471
$ws.Range("G410").Value = 628
$ws.Range("F411").Value = 224931
$ws.Range("G411").Value = 827
$ws.Range("F412").Value = 175841
$ws.Range("G412").Value = 645
$ws.Range("F413").Value = 149177
$ws.Range("F414").Value = 148615
$ws.Range("G414").Value = 564
$ws.Range("F415").Value = 305700
$ws.Range("G415").Value = 690
$ws.Range("F416").Value = 668468
$ws.Range("G416").Value = 932
$ws.Range("F417").Value = 341021
$ws.Range("G417").Value = 589
$ws.Range("F418").Value = 202139
$ws.Range("G418").Value = 702
$ws.Range("F419").Value = 148978
$ws.Range("G419").Value = 511
$ws.Range("F420").Value = 138424
$ws.Range("G420").Value = 498
$ws.Range("F421").Value = 152533
$ws.Range("G421").Value = 531
$ws.Range("F422").Value = 297119
$ws.Range("G422").Value = 645
$ws.Range("F423").Value = 434697
$ws.Range("G423").Value = 636
$ws.Range("F424").Value = 263876
$ws.Range("G424").Value = 499
$ws.Range("F425").Value = 138835
$ws.Range("G425").Value = 545
$ws.Range("F426").Value = 106724
$ws.Range("G426").Value = 386
$ws.Range("F427").Value = 90785
$ws.Range("F428").Value = 102115
$ws.Range("G428").Value = 389
$ws.Range("F429").Value = 177223
$ws.Range("G429").Value = 460
$ws.Range("F430").Value = 172863
$ws.Range("F431").Value = 170220
$ws.Range("G431").Value = 399
$ws.Range("F432").Value = 123213
$ws.Range("G432").Value = 428
$ws.Range("F433").Value = 85888
$ws.Range("G433").Value = 262
$ws.Range("F434").Value = 79537
$ws.Range("G434").Value = 280
$ws.Range("F435").Value = 81817
$ws.Range("G435").Value = 266
$ws.Range("F436").Value = 141127
$ws.Range("G436").Value = 344
$ws.Range("F437").Value = 163262
$ws.Range("G437").Value = 273
$ws.Range("F438").Value = 114820
$ws.Range("G438").Value = 243

# --- Append new row 439 (new day: 2021-05-18 / serial 44333) ---
$ws.Range("A439").Value = 44333
$ws.Range("B439").Value = 387892
$ws.Range("C439").Value = 7276
$ws.Range("D439").Value = 233
$ws.Range("E439").Value = 12248
$ws.Range("F439").Value = 70972
$ws.Range("G439").Value = 265
